# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets.
# Both sheets carry duplicate data; apply the same updates to each.

$wb = $excel.ActiveWorkbook

$updates = @{
    5  = 78
    7  = 1206
    8  = 1499
    10 = 372
    13 = 163
    15 = 103
    16 = 268
    19 = 1703
    23 = 649
    25 = 331
    26 = 4093
    30 = 1064
    31 = 131
    33 = 431
    35 = 196
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
